# "fixed optimal path for map1"
#
# Map 1 lives in B2:F6 (E = start/goal marker at B2 and F6, "*" marks the
# solved path). The previous path only ran along row 2 and column F; it
# skipped connecting down column B and across row 6 to reach the "S" cell
# at F6. Extend the "*" markers down B3:B6 and across C6:E6 to complete
# the optimal route.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3:B6").Value = "*"
$ws.Range("C6:E6").Value = "*"

# Author's saved cursor position ended up on B3 after editing the path.
$ws.Range("B3").Select()
